# Insert a new data row at row 251 (shifting existing rows 251:287 down to 252:288)
# and populate it with the new Ají / "Americana (o)" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 251; this shifts rows 251-287
# down to 252-288 and extends the used range to A1:R288.
$ws.Rows.Item(251).Insert()

# Populate the newly inserted row 251 with the new record's values.
$ws.Cells.Item(251, 1).Value = 5
$ws.Cells.Item(251, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(251, 3).Value = "Maule"
$ws.Cells.Item(251, 4).Value = 44946
$ws.Cells.Item(251, 5).Value = 7
$ws.Cells.Item(251, 6).Value = 100112021
$ws.Cells.Item(251, 7).Value = "Ají"
$ws.Cells.Item(251, 8).Value = "Americana (o)"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 100
$ws.Cells.Item(251, 11).Value = 10000
$ws.Cells.Item(251, 12).Value = 10000
$ws.Cells.Item(251, 13).Value = 10000
$ws.Cells.Item(251, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(251, 15).Value = "Región del Maule"
$ws.Cells.Item(251, 16).Value = 667
$ws.Cells.Item(251, 17).Value = 15
$ws.Cells.Item(251, 18).Value = "Hortaliza"
